$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Leading apostrophe forces Excel to store a purely-numeric-looking
    # string (e.g. "253.62") as text rather than converting it to a
    # number, matching the workbook's existing convention where every
    # "Price" cell is text-typed. Resetting the style back to Normal
    # afterwards clears the transient "quote prefix" cell style that
    # assigning a leading apostrophe would otherwise leave behind, so
    # the cell's style index is left untouched, same as before the edit.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# D/E column updates (price + volume) for most rows
Set-TextValue "D2" "42.233.94"
$ws.Range("E2").Value = "  +1.49%  "

Set-TextValue "D3" "2.173.50"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue "D5" "253.62"
$ws.Range("E5").Value = "  +6.65%  "

Set-TextValue "D6" "0.604"
$ws.Range("E6").Value = "  -0.89%  "

Set-TextValue "D7" "73.53"
$ws.Range("E7").Value = "  +1.78%  "

Set-TextValue "D9" "0.581"
$ws.Range("E9").Value = "  +0.55%  "

Set-TextValue "D10" "39.93"
$ws.Range("E10").Value = "  +0.50%  "

Set-TextValue "D11" "0.0908"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("E12").Value = "  +1.18%  "

Set-TextValue "D13" "6.76"
$ws.Range("E13").Value = "  +0.99%  "

Set-TextValue "D14" "2.500.36"
$ws.Range("E14").Value = "  +0.41%  "

Set-TextValue "D15" "14.15"
$ws.Range("E15").Value = "  -0.77%  "

Set-TextValue "D16" "2.150.64"
$ws.Range("E16").Value = "  -0.06%  "

Set-TextValue "D17" "0.764"
$ws.Range("E17").Value = "  -1.60%  "

Set-TextValue "D18" "42.116.65"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("E19").Value = "  +0.07%  "

Set-TextValue "D20" "70.49"
$ws.Range("E20").Value = "  +0.73%  "

Set-TextValue "D21" "5.83"
$ws.Range("E21").Value = "  +0.90%  "

Set-TextValue "D22" "226.52"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  -4.22%  "

Set-TextValue "D24" "2.15"
$ws.Range("E24").Value = "  +6.72%  "

$ws.Range("E25").Value = "  -0.10%  "

Set-TextValue "D26" "10.42"
$ws.Range("E26").Value = "  -2.66%  "

Set-TextValue "D27" "3.32"
$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("E28").Value = "  +2.72%  "

$ws.Range("E29").Value = "  +0.56%  "

Set-TextValue "D30" "36.79"
$ws.Range("E30").Value = "  +13.12%  "

Set-TextValue "D31" "168.50"
$ws.Range("E31").Value = "  -1.35%  "

Set-TextValue "D32" "19.98"
$ws.Range("E32").Value = "  +0.90%  "

Set-TextValue "D33" "0.0805"
$ws.Range("E33").Value = "  +4.23%  "

Set-TextValue "D34" "5.12"
$ws.Range("E34").Value = "  -4.85%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("E36").Value = "  +4.51%  "

$ws.Range("E37").Value = "  -1.32%  "

Set-TextValue "D38" "0.0328"
$ws.Range("E38").Value = "  +8.09%  "

Set-TextValue "D39" "11.85"
$ws.Range("E39").Value = "  -2.24%  "

$ws.Range("E40").Value = "  -2.50%  "

Set-TextValue "D41" "0.196"
$ws.Range("E41").Value = "  +4.01%  "

Set-TextValue "D42" "58.74"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  -4.54%  "

Set-TextValue "D44" "102.41"
$ws.Range("E44").Value = "  +5.13%  "

# Rows 45 & 46 swap: FraxShare now ranks above WOONetwork
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D45" "8.26"
$ws.Range("E45").Value = "  -1.88%  "

$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue "D46" "0.464"
$ws.Range("E46").Value = "  +15.63%  "

Set-TextValue "D47" "0.0966"
$ws.Range("E47").Value = "  +0.21%  "

Set-TextValue "D48" "2.40"
$ws.Range("E48").Value = "  +9.84%  "

Set-TextValue "D49" "1.09"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("E51").Value = "  +0.95%  "
